$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "MISSED: $find"
    }
}

function InsertRun($pos, $text) {
    # Insert $text at the collapsed point $pos, then force it to become its
    # own run (with plain/empty rPr) by toggling a character property on the
    # exact inserted span -- this runtime auto-merges runs that share
    # identical formatting, so a quick Bold on/off splits the run without
    # leaving any visible formatting behind.
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($text)
    $splitRange = $d.Range($pos, $pos + $text.Length)
    $splitRange.Bold = 1
    $splitRange.Bold = 0
    return $pos + $text.Length
}

# --- 1-6: merge the split numbered-list runs back into single runs (pure
# run-coalescing; the text is unchanged) ---------------------------------
ReplaceText "1. first" "1. first"
ReplaceText "2. The root" "2. The root"
ReplaceText "3. All leaves" "3. All leaves"
ReplaceText "4. If a node" "4. If a node"
ReplaceText "5. Every path" "5. Every path"
ReplaceText "6. the number" "6. the number"

# --- 7: "My Accomplishment" paragraph - merge the run after the <w:br/>
# with the following run (same empty rPr) -----------------------------
ReplaceText "the test. In ad" "the test. In ad"

# --- 8: "Performance of the tester" paragraph - trim the original run and
# append new content as new runs/paragraphs ------------------------------
$tail = "I found some problem, it sometimes will show the weird output in when I run the random tester, it happened because of the add or operational step proceed forward to the initialization, sometimes it works when I using ~ in front of the RedBlackTree.add().  "
$r = $d.Content
$found = $r.Find.Execute($tail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "MISSED: tail text" }
$r.Delete()

$r2 = $d.Content
$found2 = $r2.Find.Execute("but in this circumstance, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { Write-Output "MISSED: circumstance anchor" }
$r2.Collapse(0)
$pos = $r2.Start

$pos = InsertRun $pos "The random tester is not perfect."

$br = $d.Range($pos, $pos)
$br.InsertParagraphAfter()
$pos = $pos + 1

$newRuns = @(
  "The operations that wait for test is randomly picked, but ",
  "sometimes ",
  "the tester ",
  "will ",
  "choose the operation no follow the resonable order, for instance",
  ", ",
  "the ",
  "add",
  "ing",
  " ",
  "and",
  " operational step proceed forward to the initialization, ",
  "the tester will detect the error and stop running, but the program may not have errors, ",
  "sometimes it works when I using ~ in front of the RedBlackTree.add(), ",
  "it prevent the actions process before there have no initialization",
  ".  "
)
foreach ($t in $newRuns) {
    $pos = InsertRun $pos $t
}

# --- 9: "What is good about TSTL" paragraph - merge all the split runs --
ReplaceText "This random tester has its advantages" "This random tester has its advantages"

# --- 10: "Suggestion for Future work" paragraph - merge "F" + "or instance..."
ReplaceText "For instance, I realize" "For instance, I realize"

# --- "Coverage summery" heading - merge "C" + "overage " + "summery" ----
ReplaceText "Coverage summery" "Coverage summery"

# --- "We consider this part..." paragraph - merge "W" + "e consider..." -
ReplaceText "We consider this part" "We consider this part"
